$wb = $excel.ActiveWorkbook

# --- Main sheet edits ---
$main = $wb.Worksheets.Item("Main")

# C5: 50 -> -50
$main.Range("C5").Value = -50

# B2: apply date format + wrap text (new style, numFmtId 165 + wrapText)
$dateCell = $main.Range("B2")
$dateCell.NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$dateCell.WrapText = $true

# Column D: remove its currency style (column-level + per-cell D1:D5)
$main.Columns.Item(4).ClearFormats()

# Selection on Main: column C selected (as if the user clicked the column header)
$main.Columns.Item(3).Select()

# --- AccountBalance sheet edits ---
$balance = $wb.Worksheets.Item("AccountBalance")
$balance.Activate()

# Row 2 grows taller (content re-wrapped after the recalculated total)
$balance.Rows.Item(2).RowHeight = 29.25

# Selection on AccountBalance: B5
$balance.Range("B5").Select()
